$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.576.28"
$ws.Range("E2").Value = "  +1.93%  "
$ws.Range("D3").Value = "2.158.35"
$ws.Range("E3").Value = "  +3.17%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'229.69"
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("D6").Value = "'0.622"
$ws.Range("E6").Value = "  +1.09%  "
$ws.Range("D7").Value = "'63.16"
$ws.Range("E7").Value = "  +4.49%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'0.394"
$ws.Range("E9").Value = "  +2.57%  "
$ws.Range("D10").Value = "'0.0861"
$ws.Range("E10").Value = "  +2.96%  "
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("D12").Value = "'16.20"
$ws.Range("E12").Value = "  +8.11%  "
$ws.Range("D13").Value = "2.478.63"
$ws.Range("E13").Value = "  +3.16%  "
$ws.Range("D14").Value = "'22.28"
$ws.Range("E14").Value = "  +1.92%  "
$ws.Range("D15").Value = "'0.821"
$ws.Range("E15").Value = "  +3.19%  "
$ws.Range("D16").Value = "'5.57"
$ws.Range("E16").Value = "  +1.89%  "
$ws.Range("D17").Value = "2.158.80"
$ws.Range("E17").Value = "  +3.14%  "
$ws.Range("D18").Value = "39.591.93"
$ws.Range("E18").Value = "  +2.42%  "
$ws.Range("D19").Value = "'72.43"
$ws.Range("E19").Value = "  +1.12%  "
$ws.Range("D20").Value = "'6.16"
$ws.Range("D21").Value = "0.0₃0854"
$ws.Range("E21").Value = "  +2.01%  "
$ws.Range("D22").Value = "'229.18"
$ws.Range("E22").Value = "  +0.81%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "'2.41"
$ws.Range("E24").Value = "  +0.83%  "
$ws.Range("D25").Value = "'2.37"
$ws.Range("E25").Value = "  +1.51%  "
$ws.Range("D26").Value = "'9.80"
$ws.Range("E26").Value = "  +3.27%  "
$ws.Range("D27").Value = "'172.76"
$ws.Range("E27").Value = "  +1.00%  "
$ws.Range("E28").Value = "  -0.83%  "
$ws.Range("E29").Value = "  -2.16%  "
$ws.Range("D30").Value = "'19.66"
$ws.Range("E30").Value = "  +2.51%  "
$ws.Range("E31").Value = "  +8.58%  "
$ws.Range("E32").Value = "  +1.46%  "
$ws.Range("D33").Value = "'4.64"
$ws.Range("E33").Value = "  +3.07%  "
$ws.Range("D34").Value = "'4.83"
$ws.Range("E34").Value = "  +2.86%  "
$ws.Range("D35").Value = "'7.17"
$ws.Range("E35").Value = "  +11.35%  "
$ws.Range("E36").Value = "  +1.98%  "
$ws.Range("D37").Value = "'2.43"
$ws.Range("E37").Value = "  +2.19%  "
$ws.Range("D38").Value = "'3.58"
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("D39").Value = "'0.999"
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("D40").Value = "'18.25"
$ws.Range("E40").Value = "  +0.48%  "
$ws.Range("D41").Value = "'0.0233"
$ws.Range("E41").Value = "  +3.67%  "
$ws.Range("D42").Value = "'103.40"
$ws.Range("E42").Value = "  +2.48%  "
$ws.Range("D43").Value = "1.533.77"
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("E44").Value = "  +5.77%  "
$ws.Range("D45").Value = "'0.0929"
$ws.Range("E45").Value = "  +0.77%  "
$ws.Range("E46").Value = "  +7.12%  "
$ws.Range("E47").Value = "  -0.36%  "
$ws.Range("D48").Value = "'7.76"
$ws.Range("E48").Value = "  +1.68%  "
$ws.Range("D49").Value = "'4.19"
$ws.Range("E49").Value = "  +2.05%  "
$ws.Range("D50").Value = "2.362.30"
$ws.Range("E50").Value = "  +3.20%  "
$ws.Range("E51").Value = "  +0.06%  "
